# Regenerate the IFRS financial figures for DB금융투자 (company_list sheet).
# Rows 2-6: numeric figures replaced with corrected values; the "FCF" (U) column
# is dropped from each row (shifted data no longer has an FCF figure separate from CAPEX).
# Rows 7-9: these annual columns are blanked out (data only kept through row 6).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update financial figures
$ws.Range("D2").Value = 10861
$ws.Range("E2").Value = 213
$ws.Range("F2").Value = 213
$ws.Range("G2").Value = 227
$ws.Range("H2").Value = 163
$ws.Range("I2").Value = 149
$ws.Range("J2").Value = 14
$ws.Range("K2").Value = 64903
$ws.Range("L2").Value = 57953
$ws.Range("M2").Value = 6949
$ws.Range("N2").Value = 6110
$ws.Range("O2").Value = 839
$ws.Range("P2").Value = 2122
$ws.Range("Q2").Value = -2611
$ws.Range("R2").Value = 883
$ws.Range("S2").Value = 3041
$ws.Range("T2").Value = 18
$ws.Range("V2").Value = 8392
$ws.Range("W2").Value = 1.97
$ws.Range("X2").Value = 1.5
$ws.Range("Y2").Value = 2.47
$ws.Range("Z2").Value = 0.24
$ws.Range("AA2").Value = 833.9299999999999
$ws.Range("AB2").Value = 231
$ws.Range("AC2").Value = 352
$ws.Range("AD2").Value = 10.54
$ws.Range("AE2").Value = 14749
$ws.Range("AF2").Value = 0.25
$ws.Range("AG2").Value = 100
$ws.Range("AH2").Value = 2.7
$ws.Range("AI2").Value = 27.73
$ws.Range("AJ2").Value = 42446389
# Row 2: remove obsolete column(s)
$ws.Range("U2").ClearContents()

# Row 3: update financial figures
$ws.Range("D3").Value = 13440
$ws.Range("E3").Value = 104
$ws.Range("F3").Value = 104
$ws.Range("G3").Value = -72
$ws.Range("H3").Value = -85
$ws.Range("I3").Value = -75
$ws.Range("J3").Value = -10
$ws.Range("K3").Value = 66392
$ws.Range("L3").Value = 59541
$ws.Range("M3").Value = 6850
$ws.Range("N3").Value = 6027
$ws.Range("O3").Value = 823
$ws.Range("P3").Value = 2122
$ws.Range("Q3").Value = -169
$ws.Range("R3").Value = -101
$ws.Range("S3").Value = -239
$ws.Range("T3").Value = 68
$ws.Range("V3").Value = 12619
$ws.Range("W3").Value = 0.77
$ws.Range("X3").Value = -0.63
$ws.Range("Y3").Value = -1.23
$ws.Range("Z3").Value = -0.11
$ws.Range("AA3").Value = 869.1799999999999
$ws.Range("AB3").Value = 226.33
$ws.Range("AC3").Value = -176
$ws.Range("AD3").Value = -22
$ws.Range("AE3").Value = 14547
$ws.Range("AF3").Value = 0.27
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 42446389
# Row 3: remove obsolete column(s)
$ws.Range("U3").ClearContents()

# Row 4: update financial figures
$ws.Range("D4").Value = 13254
$ws.Range("E4").Value = 98
$ws.Range("F4").Value = 98
$ws.Range("G4").Value = 125
$ws.Range("H4").Value = 64
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 63
$ws.Range("K4").Value = 63356
$ws.Range("L4").Value = 56396
$ws.Range("M4").Value = 6960
$ws.Range("N4").Value = 6081
$ws.Range("O4").Value = 879
$ws.Range("P4").Value = 2122
$ws.Range("Q4").Value = 3077
$ws.Range("R4").Value = 35
$ws.Range("S4").Value = -1445
$ws.Range("T4").Value = 29
$ws.Range("V4").Value = 11345
$ws.Range("W4").Value = 0.74
$ws.Range("X4").Value = 0.49
$ws.Range("Y4").Value = 0.02
$ws.Range("Z4").Value = 0
$ws.Range("AA4").Value = 810.25
$ws.Range("AB4").Value = 231.51
$ws.Range("AC4").Value = 2
$ws.Range("AD4").Value = 1308.55
$ws.Range("AE4").Value = 14678
$ws.Range("AF4").Value = 0.22
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 42446389
# Row 4: remove obsolete column(s)
$ws.Range("U4").ClearContents()

# Row 5: update financial figures
$ws.Range("D5").Value = 9099
$ws.Range("E5").Value = 224
$ws.Range("F5").Value = 224
$ws.Range("G5").Value = 251
$ws.Range("H5").Value = 154
$ws.Range("I5").Value = 44
$ws.Range("J5").Value = 109
$ws.Range("K5").Value = 69983
$ws.Range("L5").Value = 62907
$ws.Range("M5").Value = 7077
$ws.Range("N5").Value = 6096
$ws.Range("O5").Value = 981
$ws.Range("P5").Value = 2122
$ws.Range("Q5").Value = -1315
$ws.Range("R5").Value = 68
$ws.Range("S5").Value = 1395
$ws.Range("T5").Value = 18
$ws.Range("V5").Value = 12883
$ws.Range("W5").Value = 2.46
$ws.Range("X5").Value = 1.69
$ws.Range("Y5").Value = 0.73
$ws.Range("Z5").Value = 0.07000000000000001
$ws.Range("AA5").Value = 888.92
$ws.Range("AB5").Value = 237
$ws.Range("AC5").Value = 104
$ws.Range("AD5").Value = 38.34
$ws.Range("AE5").Value = 14713
$ws.Range("AF5").Value = 0.27
$ws.Range("AG5").Value = 100
$ws.Range("AH5").Value = 2.5
$ws.Range("AI5").Value = 93.44
$ws.Range("AJ5").Value = 42446389
# Row 5: remove obsolete column(s)
$ws.Range("U5").ClearContents()

# Row 6: update financial figures
$ws.Range("D6").Value = 8487
$ws.Range("E6").Value = 864
$ws.Range("F6").Value = 864
$ws.Range("G6").Value = 870
$ws.Range("H6").Value = 631
$ws.Range("I6").Value = 543
$ws.Range("K6").Value = 70432
$ws.Range("L6").Value = 62823
$ws.Range("M6").Value = 7609
$ws.Range("N6").Value = 6577
$ws.Range("P6").Value = 2122
$ws.Range("Q6").Value = 1688
$ws.Range("R6").Value = -24
$ws.Range("S6").Value = -1598
$ws.Range("T6").Value = 25
$ws.Range("V6").Value = 12041
$ws.Range("W6").Value = 10.18
$ws.Range("X6").Value = 7.44
$ws.Range("Y6").Value = 8.58
$ws.Range("Z6").Value = 0.77
$ws.Range("AA6").Value = 825.6799999999999
$ws.Range("AB6").Value = 262.06
$ws.Range("AC6").Value = 1280
$ws.Range("AD6").Value = 3.47
$ws.Range("AE6").Value = 15875
$ws.Range("AF6").Value = 0.28
$ws.Range("AG6").Value = 250
$ws.Range("AH6").Value = 5.62
$ws.Range("AI6").Value = 19.06
$ws.Range("AJ6").Value = 42446389
# Row 6: remove obsolete column(s)
$ws.Range("U6").ClearContents()

# Row 7: clear all financial figures (data series ends at row 6)
$ws.Range("D7:AJ7").ClearContents()

# Row 8: clear all financial figures (data series ends at row 6)
$ws.Range("D8:AJ8").ClearContents()

# Row 9: clear all financial figures (data series ends at row 6)
$ws.Range("D9:AJ9").ClearContents()
